# Generate Report for Handoff
# This script updates the localization-status workbook so that the row
# previously describing "eaebf39c-...md" (handed back, in sync) now sits
# in row 2 and the row describing "39b07019-...md" moves to row 3 with an
# updated status ("Ready for handoff") and refreshed handoff timestamps.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A2").Value = "eaebf39c-c09c-44b7-9f1e-18c264aef4fe.md"
$ws1.Range("B2").Value = "Handed back: in sync with en-US"
$ws1.Range("C2").Value = "Handed back: in sync with en-US"
$ws1.Range("D2").Value = "2016-48-13 06:48:01"

$ws1.Range("A3").Value = "39b07019-896a-4d16-842b-bb42829f0703.md"
$ws1.Range("B3").Value = "Ready for handoff"
$ws1.Range("C3").Value = "Ready for handoff"
$ws1.Range("D3").Value = "2016-50-13 06:50:43"

$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/67ce2f41856542034e1c7a32912c1dbbf1a41f75/e2e/39b07019-896a-4d16-842b-bb42829f0703.md", [Type]::Missing, [Type]::Missing, "eaebf39c-c09c-44b7-9f1e-18c264aef4fe.md")
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/67ce2f41856542034e1c7a32912c1dbbf1a41f75/e2e/eaebf39c-c09c-44b7-9f1e-18c264aef4fe.md", [Type]::Missing, [Type]::Missing, "39b07019-896a-4d16-842b-bb42829f0703.md")

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A2").Value = "eaebf39c-c09c-44b7-9f1e-18c264aef4fe.md"
$ws2.Range("B2").Value = ".md"
$ws2.Range("C2").Value = "Handed back: in sync with en-US"
$ws2.Range("D2").Value = "eaebf39c-c09c-44b7-9f1e-18c264aef4fe.ca39043310146e438fc77fd927f1989d85b3784b.zh-cn.xlf"
$ws2.Range("E2").Value = "2016-03-13 06:46:32"
$ws2.Range("F2").Value = "eaebf39c-c09c-44b7-9f1e-18c264aef4fe.md"
$ws2.Range("G2").Value = "eaebf39c-c09c-44b7-9f1e-18c264aef4fe.ca39043310146e438fc77fd927f1989d85b3784b.zh-cn.xlf"
$ws2.Range("H2").Value = "2016-03-13 06:50:14"
$ws2.Range("I2").Value = "Include"

$ws2.Range("A3").Value = "39b07019-896a-4d16-842b-bb42829f0703.md"
$ws2.Range("B3").Value = ".md"
$ws2.Range("C3").Value = "Ready for handoff"
$ws2.Range("D3").Value = "39b07019-896a-4d16-842b-bb42829f0703.14ddb031785dcc11513facc8483568f5197f3d09.zh-cn.xlf"
$ws2.Range("E3").Value = "2016-03-13 06:50:39"
$ws2.Range("F3").Value = "39b07019-896a-4d16-842b-bb42829f0703.md"
$ws2.Range("G3").Value = "39b07019-896a-4d16-842b-bb42829f0703.14ddb031785dcc11513facc8483568f5197f3d09.zh-cn.xlf"
$ws2.Range("H3").Value = "2016-03-13 06:50:14"
$ws2.Range("I3").Value = "Include"

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/67ce2f41856542034e1c7a32912c1dbbf1a41f75/e2e/39b07019-896a-4d16-842b-bb42829f0703.md", [Type]::Missing, [Type]::Missing, "eaebf39c-c09c-44b7-9f1e-18c264aef4fe.md")
$ws2.Hyperlinks.Add($ws2.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/67ce2f41856542034e1c7a32912c1dbbf1a41f75/e2e/39b07019-896a-4d16-842b-bb42829f0703.md", [Type]::Missing, [Type]::Missing, ".md")
$ws2.Hyperlinks.Add($ws2.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fa66c6db28e2bbb0ef480c9bfc26f1bdc51bb086/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/39b07019-896a-4d16-842b-bb42829f0703.14ddb031785dcc11513facc8483568f5197f3d09.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "eaebf39c-c09c-44b7-9f1e-18c264aef4fe.ca39043310146e438fc77fd927f1989d85b3784b.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/0e51df6a84d113a7a43484d1cf13e3e52716eac9/e2e/39b07019-896a-4d16-842b-bb42829f0703.md", [Type]::Missing, [Type]::Missing, "eaebf39c-c09c-44b7-9f1e-18c264aef4fe.md")
$ws2.Hyperlinks.Add($ws2.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/b045de3e33a694e56d697bcfa7777fec0cc03097/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/39b07019-896a-4d16-842b-bb42829f0703.14ddb031785dcc11513facc8483568f5197f3d09.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "eaebf39c-c09c-44b7-9f1e-18c264aef4fe.ca39043310146e438fc77fd927f1989d85b3784b.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/67ce2f41856542034e1c7a32912c1dbbf1a41f75/e2e/eaebf39c-c09c-44b7-9f1e-18c264aef4fe.md", [Type]::Missing, [Type]::Missing, "39b07019-896a-4d16-842b-bb42829f0703.md")
$ws2.Hyperlinks.Add($ws2.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/67ce2f41856542034e1c7a32912c1dbbf1a41f75/e2e/eaebf39c-c09c-44b7-9f1e-18c264aef4fe.md", [Type]::Missing, [Type]::Missing, ".md")
$ws2.Hyperlinks.Add($ws2.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fa66c6db28e2bbb0ef480c9bfc26f1bdc51bb086/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/eaebf39c-c09c-44b7-9f1e-18c264aef4fe.ca39043310146e438fc77fd927f1989d85b3784b.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "39b07019-896a-4d16-842b-bb42829f0703.14ddb031785dcc11513facc8483568f5197f3d09.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/0e51df6a84d113a7a43484d1cf13e3e52716eac9/e2e/eaebf39c-c09c-44b7-9f1e-18c264aef4fe.md", [Type]::Missing, [Type]::Missing, "39b07019-896a-4d16-842b-bb42829f0703.md")
$ws2.Hyperlinks.Add($ws2.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/b045de3e33a694e56d697bcfa7777fec0cc03097/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/eaebf39c-c09c-44b7-9f1e-18c264aef4fe.ca39043310146e438fc77fd927f1989d85b3784b.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "39b07019-896a-4d16-842b-bb42829f0703.14ddb031785dcc11513facc8483568f5197f3d09.zh-cn.xlf")

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A2").Value = "eaebf39c-c09c-44b7-9f1e-18c264aef4fe.md"
$ws3.Range("B2").Value = ".md"
$ws3.Range("C2").Value = "Handed back: in sync with en-US"
$ws3.Range("D2").Value = "eaebf39c-c09c-44b7-9f1e-18c264aef4fe.ca39043310146e438fc77fd927f1989d85b3784b.de-de.xlf"
$ws3.Range("E2").Value = "2016-03-13 06:48:01"
$ws3.Range("F2").Value = "eaebf39c-c09c-44b7-9f1e-18c264aef4fe.md"
$ws3.Range("G2").Value = "eaebf39c-c09c-44b7-9f1e-18c264aef4fe.ca39043310146e438fc77fd927f1989d85b3784b.de-de.xlf"
$ws3.Range("H2").Value = "2016-03-13 06:50:21"
$ws3.Range("I2").Value = "Include"

$ws3.Range("A3").Value = "39b07019-896a-4d16-842b-bb42829f0703.md"
$ws3.Range("B3").Value = ".md"
$ws3.Range("C3").Value = "Ready for handoff"
$ws3.Range("D3").Value = "39b07019-896a-4d16-842b-bb42829f0703.14ddb031785dcc11513facc8483568f5197f3d09.de-de.xlf"
$ws3.Range("E3").Value = "2016-03-13 06:50:43"
$ws3.Range("F3").Value = "39b07019-896a-4d16-842b-bb42829f0703.md"
$ws3.Range("G3").Value = "39b07019-896a-4d16-842b-bb42829f0703.14ddb031785dcc11513facc8483568f5197f3d09.de-de.xlf"
$ws3.Range("H3").Value = "2016-03-13 06:50:21"
$ws3.Range("I3").Value = "Include"

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/67ce2f41856542034e1c7a32912c1dbbf1a41f75/e2e/39b07019-896a-4d16-842b-bb42829f0703.md", [Type]::Missing, [Type]::Missing, "eaebf39c-c09c-44b7-9f1e-18c264aef4fe.md")
$ws3.Hyperlinks.Add($ws3.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/67ce2f41856542034e1c7a32912c1dbbf1a41f75/e2e/39b07019-896a-4d16-842b-bb42829f0703.md", [Type]::Missing, [Type]::Missing, ".md")
$ws3.Hyperlinks.Add($ws3.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/01aa9f03f0dbd5ab36c898cd118e5fc6c73d2b01/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/39b07019-896a-4d16-842b-bb42829f0703.14ddb031785dcc11513facc8483568f5197f3d09.de-de.xlf", [Type]::Missing, [Type]::Missing, "eaebf39c-c09c-44b7-9f1e-18c264aef4fe.ca39043310146e438fc77fd927f1989d85b3784b.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/4665ac3c2f0d624e0ca2b515bd68029ddfaeb16a/e2e/39b07019-896a-4d16-842b-bb42829f0703.md", [Type]::Missing, [Type]::Missing, "eaebf39c-c09c-44b7-9f1e-18c264aef4fe.md")
$ws3.Hyperlinks.Add($ws3.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/474830c5cf9b4ed106186dca00bc03643d35a40b/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/39b07019-896a-4d16-842b-bb42829f0703.14ddb031785dcc11513facc8483568f5197f3d09.de-de.xlf", [Type]::Missing, [Type]::Missing, "eaebf39c-c09c-44b7-9f1e-18c264aef4fe.ca39043310146e438fc77fd927f1989d85b3784b.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/67ce2f41856542034e1c7a32912c1dbbf1a41f75/e2e/eaebf39c-c09c-44b7-9f1e-18c264aef4fe.md", [Type]::Missing, [Type]::Missing, "39b07019-896a-4d16-842b-bb42829f0703.md")
$ws3.Hyperlinks.Add($ws3.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/67ce2f41856542034e1c7a32912c1dbbf1a41f75/e2e/eaebf39c-c09c-44b7-9f1e-18c264aef4fe.md", [Type]::Missing, [Type]::Missing, ".md")
$ws3.Hyperlinks.Add($ws3.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/01aa9f03f0dbd5ab36c898cd118e5fc6c73d2b01/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/eaebf39c-c09c-44b7-9f1e-18c264aef4fe.ca39043310146e438fc77fd927f1989d85b3784b.de-de.xlf", [Type]::Missing, [Type]::Missing, "39b07019-896a-4d16-842b-bb42829f0703.14ddb031785dcc11513facc8483568f5197f3d09.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/4665ac3c2f0d624e0ca2b515bd68029ddfaeb16a/e2e/eaebf39c-c09c-44b7-9f1e-18c264aef4fe.md", [Type]::Missing, [Type]::Missing, "39b07019-896a-4d16-842b-bb42829f0703.md")
$ws3.Hyperlinks.Add($ws3.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/474830c5cf9b4ed106186dca00bc03643d35a40b/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/eaebf39c-c09c-44b7-9f1e-18c264aef4fe.ca39043310146e438fc77fd927f1989d85b3784b.de-de.xlf", [Type]::Missing, [Type]::Missing, "39b07019-896a-4d16-842b-bb42829f0703.14ddb031785dcc11513facc8483568f5197f3d09.de-de.xlf")
